$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current "extra track" notes from column D before moving them
$d5 = $ws.Range("D5").Value2
$d7 = $ws.Range("D7").Value2
$d9 = $ws.Range("D9").Value2
$d11 = $ws.Range("D11").Value2

# Shift the notes up (close the gaps)
$ws.Range("D3").Value = $d5
$ws.Range("D4").Value = $d7
$ws.Range("D5").Value = $d9
$ws.Range("D6").Value = $d11

$ws.Range("D7").Clear()
$ws.Range("D11").Clear()

# Row 46: replace the "-> Inescapable Fate" note (col C) with a renamed
# version of the track name (col B)
$ws.Range("C46").Clear()
$ws.Range("B46").Value = "FE7 Inescapable Fate"

# Add a new "Tracks Remaining" note
$ws.Range("D9").Value = "Tracks Remaining: 12"

# Move the active selection to E28
$ws.Range("E28").Select() | Out-Null
